$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Prepare destination rows by copying formats from existing styled rows ---
# Capture the old "total" style (row 18) for the new total row 21 BEFORE row 18's
# own style gets overwritten below.
$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B21:J21").PasteSpecial(-4122) | Out-Null
# Rows 16-20 all use the plain "normal" row style (like old rows 16/17); row 18
# switches from the old "total" style to this plain style.
$ws.Range("B16:J17").Copy() | Out-Null
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null
$ws.Range("B19:J20").PasteSpecial(-4122) | Out-Null
# Footer rows 26-27 get the same style as old footer rows 23-24
$ws.Range("B23:J24").Copy() | Out-Null
$ws.Range("B26:J27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2) Remove the old footer rows (23-24): unmerge then clear ---
$ws.Range("B23:C23").UnMerge() | Out-Null
$ws.Range("H23:J23").UnMerge() | Out-Null
$ws.Range("B24:C24").UnMerge() | Out-Null
$ws.Range("H24:J24").UnMerge() | Out-Null
$ws.Range("B23:J24").Clear() | Out-Null

# --- 3) Re-create merges for the new footer row positions ---
$ws.Range("B26:C26").Merge() | Out-Null
$ws.Range("H26:J26").Merge() | Out-Null
$ws.Range("B27:C27").Merge() | Out-Null
$ws.Range("H27:J27").Merge() | Out-Null

# --- 4) Header values ---
$ws.Range("E11").Value = 299494
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 4

# --- 5) Worker table rows 16-21 ---
# Row 16
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73183038"
$ws.Range("D16").Value = "JULIO MANUEL MORON BATISTA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Row 17
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73183038"
$ws.Range("D17").Value = "JULIO MANUEL MORON BATISTA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73183038"
$ws.Range("D18").Value = "JULIO MANUEL MORON BATISTA"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19 (new)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047459861"
$ws.Range("D19").Value = "JHON JAIRO BATISTA MURILLO"
$ws.Range("E19").Value = "2503"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

# Row 20 (new)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "9237009"
$ws.Range("D20").Value = "ALBEIRO BATISTA OTERO"
$ws.Range("E20").Value = "2503"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

# Row 21 (totals / highlighted row, new position)
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1193596395"
$ws.Range("D21").Value = "MAIRON SANTIAGO AREVALO BATISTA"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 24674
$ws.Range("G21").Value = 1423500

# --- 6) Footer text (rows 26-27, same text as before, shifted) ---
$ws.Range("B26").Value = "___________________________________"
$ws.Range("H26").Value = "___________________________________"
$ws.Range("B27").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H27").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# --- 7) Column D width (auto-fit to longest new name) ---
$ws.Columns.Item(4).ColumnWidth = 34.3
